# Updated cryptos list on Mon Dec 25 07:49:38 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each assignment below uses a leading apostrophe inside the string so that
# numeric-looking text (prices such as "6.20", "1.00", "43.279.91") is stored
# by Excel as literal text, exactly preserving the source formatting, instead
# of being auto-converted into a number (which would drop trailing zeros).

$ws.Range("D2").Value = '''43.279.91'
$ws.Range("E2").Value = '''  -1.07%  '
$ws.Range("D3").Value = '''2.283.60'
$ws.Range("E3").Value = '''  -0.35%  '
$ws.Range("E4").Value = '''  -0.02%  '
$ws.Range("D5").Value = '''113.68'
$ws.Range("E5").Value = '''  -0.49%  '
$ws.Range("D6").Value = '''265.45'
$ws.Range("E6").Value = '''  -1.21%  '
$ws.Range("D7").Value = '''0.617'
$ws.Range("E7").Value = '''  -0.96%  '
$ws.Range("E8").Value = '''  +0.01%  '
$ws.Range("E9").Value = '''  -1.49%  '
$ws.Range("D10").Value = '''47.58'
$ws.Range("E10").Value = '''  -1.01%  '
$ws.Range("D11").Value = '''0.0935'
$ws.Range("E11").Value = '''  -0.37%  '
$ws.Range("D12").Value = '''9.29'
$ws.Range("E12").Value = '''  +7.40%  '
$ws.Range("E13").Value = '''  +1.08%  '
$ws.Range("E14").Value = '''  -0.61%  '
$ws.Range("D15").Value = '''2.611.24'
$ws.Range("E15").Value = '''  -0.87%  '
$ws.Range("D16").Value = '''0.868'
$ws.Range("E16").Value = '''  +2.30%  '
$ws.Range("D17").Value = '''2.288.49'
$ws.Range("E17").Value = '''  -0.05%  '
$ws.Range("D18").Value = '''43.244.62'
$ws.Range("E18").Value = '''  -0.84%  '
$ws.Range("E19").Value = '''  -1.02%  '
$ws.Range("E20").Value = '''  +3.52%  '
$ws.Range("D21").Value = '''71.98'
$ws.Range("E21").Value = '''  -0.67%  '
$ws.Range("D22").Value = '''2.49'
$ws.Range("E22").Value = '''  -0.73%  '
$ws.Range("D23").Value = '''233.92'
$ws.Range("E23").Value = '''  +0.24%  '
$ws.Range("D24").Value = '''9.62'
$ws.Range("E24").Value = '''  +1.08%  '
$ws.Range("E25").Value = '''  +0.59%  '
$ws.Range("E26").Value = '''  +1.68%  '
$ws.Range("D27").Value = '''11.43'
$ws.Range("E27").Value = '''  -1.26%  '
$ws.Range("D28").Value = '''3.99'
$ws.Range("E28").Value = '''  +0.05%  '
$ws.Range("D29").Value = '''41.11'
$ws.Range("E29").Value = '''  -2.23%  '
$ws.Range("D30").Value = '''3.33'
$ws.Range("E30").Value = '''  -2.21%  '
$ws.Range("E31").Value = '''  -0.78%  '
$ws.Range("D32").Value = '''173.84'
$ws.Range("E32").Value = '''  -1.33%  '
$ws.Range("D33").Value = '''21.67'
$ws.Range("E33").Value = '''  +0.30%  '
$ws.Range("D34").Value = '''0.0907'
$ws.Range("E34").Value = '''  -2.62%  '
$ws.Range("E35").Value = '''  +3.41%  '
$ws.Range("D36").Value = '''0.128'
$ws.Range("E36").Value = '''  +0.94%  '
$ws.Range("D37").Value = '''4.66'
$ws.Range("E37").Value = '''  -0.84%  '
$ws.Range("D38").Value = '''0.0368'
$ws.Range("E38").Value = '''  +2.68%  '
$ws.Range("D39").Value = '''3.93'
$ws.Range("E39").Value = '''  +2.60%  '
$ws.Range("E40").Value = '''  -5.38%  '
$ws.Range("D41").Value = '''2.61'
$ws.Range("E41").Value = '''  +8.87%  '
$ws.Range("D42").Value = '''76.76'
$ws.Range("E42").Value = '''  +3.25%  '
$ws.Range("E43").Value = '''  +1.90%  '
$ws.Range("E44").Value = '''  -1.88%  '
$ws.Range("D45").Value = '''6.20'
$ws.Range("E45").Value = '''  +3.41%  '
$ws.Range("D46").Value = '''0.999'
$ws.Range("E46").Value = '''  -0.26%  '
$ws.Range("E47").Value = '''  -3.64%  '
$ws.Range("D48").Value = '''103.76'
$ws.Range("E48").Value = '''  +1.03%  '
$ws.Range("D49").Value = '''8.57'
$ws.Range("E49").Value = '''  -2.18%  '
$ws.Range("B50").Value = '''Cronos'
$ws.Range("C50").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0997'
$ws.Range("E50").Value = '''  -0.80%  '
$ws.Range("B51").Value = '''TrustWalletToken'
$ws.Range("C51").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '''1.25'
$ws.Range("E51").Value = '''  +1.43%  '
